$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "61.882.58"
$ws.Cells.Item(2, 5).Value = "  -2.24%  "
$ws.Cells.Item(3, 4).Value = "2.575.86"
$ws.Cells.Item(3, 5).Value = "  -3.92%  "
$ws.Cells.Item(4, 5).Value = "  -0.02%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "550.76"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -0.42%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "154.85"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -1.87%  "
$ws.Cells.Item(7, 5).Value = "  -0.03%  "
$ws.Cells.Item(8, 5).Value = "  +1.69%  "
$ws.Cells.Item(9, 5).Value = "  -1.09%  "
$ws.Cells.Item(10, 5).Value = "  -1.07%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "5.50"
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +3.33%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.366"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  -0.55%  "
$ws.Cells.Item(13, 4).Value = "3.030.22"
$ws.Cells.Item(13, 5).Value = "  -3.90%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "25.42"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -3.15%  "
$ws.Cells.Item(15, 4).Value = "61.807.62"
$ws.Cells.Item(15, 5).Value = "  -2.16%  "
$ws.Cells.Item(16, 5).Value = "  -0.30%  "
$ws.Cells.Item(17, 4).Value = "2.578.56"
$ws.Cells.Item(17, 5).Value = "  -3.77%  "
$ws.Cells.Item(18, 5).Value = "  -3.20%  "
$ws.Cells.Item(19, 5).Value = "  -0.40%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "337.99"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  -1.61%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "6.02"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "1.00"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  +0.45%  "
$ws.Cells.Item(23, 5).Value = "  -3.07%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "63.49"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -0.27%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "0.167"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -0.70%  "
$ws.Cells.Item(26, 5).Value = "  -0.31%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "8.13"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "7.23"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +3.67%  "
$ws.Cells.Item(29, 4).Value = "0.0₃0833"
$ws.Cells.Item(29, 5).Value = "  -2.22%  "
$ws.Cells.Item(30, 5).Value = "  -0.03%  "
$ws.Cells.Item(31, 5).Value = "  -2.23%  "
$ws.Cells.Item(32, 5).Value = "  -1.77%  "
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "4.87"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  +1.80%  "
$ws.Cells.Item(34, 5).Value = "  +0.05%  "
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "19.15"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -1.78%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "1.41"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  -1.27%  "
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "1.78"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +0.16%  "
$ws.Cells.Item(38, 2).Value = "RenderToken"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "6.03"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -0.80%  "
$ws.Cells.Item(39, 2).Value = "Bittensor"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "325.64"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -4.14%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "0.902"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -4.23%  "
$ws.Cells.Item(41, 5).Value = "  +0.35%  "
$ws.Cells.Item(42, 5).Value = "  -1.61%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "20.58"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -0.65%  "
$ws.Cells.Item(44, 5).Value = "  -0.03%  "
$ws.Cells.Item(45, 5).Value = "  -2.00%  "
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "10.94"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -1.07%  "
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.0545"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  -2.68%  "
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "0.0965"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -0.48%  "
$ws.Cells.Item(49, 5).Value = "  -4.01%  "
$ws.Cells.Item(50, 5).Value = "  -1.30%  "
$ws.Cells.Item(51, 4).Value = "2.048.56"
$ws.Cells.Item(51, 5).Value = "  -1.90%  "
